# Update daily COVID stats worksheet: po 11. 01. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrections to existing rows (columns H = AgTests, I = AgPosit)
$ws.Range("H293").Value = 82469

$ws.Range("H294").Value = 91708
$ws.Range("I294").Value = 5088

$ws.Range("H299").Value = 64720
$ws.Range("I299").Value = 6798

$ws.Range("H300").Value = 70631
$ws.Range("I300").Value = 6919

$ws.Range("H301").Value = 69560
$ws.Range("I301").Value = 5526

$ws.Range("H302").Value = 72518
$ws.Range("I302").Value = 5281

$ws.Range("H305").Value = 3625
$ws.Range("I305").Value = 302

$ws.Range("H306").Value = 69879
$ws.Range("I306").Value = 7112

$ws.Range("H307").Value = 72032
$ws.Range("I307").Value = 6216

$ws.Range("H309").Value = 56648
$ws.Range("I309").Value = 3915

# New rows for 2021-01-08, 2021-01-09, 2021-01-10
$ws.Range("A310").Value = 44204
$ws.Range("B310").Value = 205236
$ws.Range("C310").Value = 144612
$ws.Range("D310").Value = 57788
$ws.Range("E310").Value = 16368
$ws.Range("F310").Value = 4072
$ws.Range("G310").Value = 2836
$ws.Range("H310").Value = 84945
$ws.Range("I310").Value = 5131

$ws.Range("A311").Value = 44205
$ws.Range("B311").Value = 208209
$ws.Range("C311").Value = 147275
$ws.Range("D311").Value = 58016
$ws.Range("E311").Value = 12844
$ws.Range("F311").Value = 2973
$ws.Range("G311").Value = 2918
$ws.Range("H311").Value = 22875
$ws.Range("I311").Value = 1117

$ws.Range("A312").Value = 44206
$ws.Range("B312").Value = 209069
$ws.Range("C312").Value = 150239
$ws.Range("D312").Value = 55823
$ws.Range("E312").Value = 4328
$ws.Range("F312").Value = 860
$ws.Range("G312").Value = 3007
$ws.Range("H312").Value = 31377
$ws.Range("I312").Value = 1012

# Apply same date number format (s="2") as column A for the new rows, matching existing column A formatting
$ws.Range("A310:A312").NumberFormat = "yyyy-mm-dd"
